$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 11) so the table shrinks from 11 to 10 rows
$ws.Rows.Item(11).Delete()

# New data values (rows 2-10), columns A-H
$data = @(
    @(7, "2025-04-23", 400, "V V REFEICOES LTDA", "000029", "ESPONJA MULTIUSO JEITOSA", 8168, $false),
    @(1, "2025-04-24", 150, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000088", "VASSOURA PIACAVA 20 FUROS", 185, $false),
    @(4, "2025-04-24", 300, "MUSASHI DA AMAZONIA LTDA", "000842", "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND", 636, $false),
    @(0, "2025-04-28", 250, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000098", "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM", -5, $false),
    @(2, "2025-04-28", 60, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000347", "INSETICIDA BUZZOFF AEROSSOL 300ML", 2, $true),
    @(5, "2025-04-28", 70, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000349", "DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO", 356, $true),
    @(6, "2025-04-28", 20, "LUCAS CLIENTE NOVO", "000158", "AZULIM LIMPA CERAMICAS E AZULEJOS LAVANDA 5L 1:15 START", 0, $true),
    @(8, "2025-04-28", 250, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000779", "PEDRA SANITARIA NAFT PLUS FLORAL 25G", 184, $false),
    @(3, "2025-04-30", 250, "RH MULTI SERVICOS ADMINISTRATIVOS S.A", "000041", "LUVAS DESCARTAVEIS C/ 100 UND", 1063, $false)
)

# Columns B and E hold text-like values (dates written as plain strings, and
# zero-padded numeric codes) that Excel would otherwise auto-convert to a
# date serial / number. Mark them as Text first, then restore the original
# (default) cell style after assigning the values so no stray style index
# is introduced.
$ws.Range("B2:B10").NumberFormat = "@"
$ws.Range("E2:E10").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $rowIndex++
}

$ws.Range("B2:B10").Style = "Normal"
$ws.Range("E2:E10").Style = "Normal"
